$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/border/centered
# header style already used by the other header cells (e.g. H1).
$headerRng = $ws.Range("I1:J1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108   # xlCenter
$headerRng.VerticalAlignment = -4160     # xlTop
$headerRng.Borders.LineStyle = 1         # xlContinuous
$headerRng.Borders.Weight = 2            # xlThin

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J for rows 2-4
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 5
